# Automatische test-sync: 2025-06-23 18:41:50
# Adds the new "Beschadigd product ontvangen" mail-log entry as row 19 on
# the "Logs" sheet, extends the D/G conditional-formatting ranges to cover
# it, and swaps the Dashboard's "Bestelling / Levering" / "Retour /
# Terugbetaling" rows (with the refreshed counts) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- New row 19 on the "Logs" sheet ---------------------------------------
$ws.Cells.Item(19, 1).Value = "Beschadigd product ontvangen"
$ws.Cells.Item(19, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(19, 3).Value = "Het product dat ik heb ontvangen is beschadigd aangekomen."
$ws.Cells.Item(19, 4).Value = "Retour / Terugbetaling"
$ws.Cells.Item(19, 5).Value = "Beste klant,
Bedankt voor uw bericht en onze excuses voor het ongemak dat u heeft ervaren. Om u verder te kunnen helpen, hebben we wat aanvullende informatie nodig. Kunt u alstublieft een foto van de beschadigde product(en) meesturen, samen met uw ordernummer? Op die manier kunnen we het probleem beter onderzoeken en zo snel mogelijk een passende oplossing bieden.
Wij kijken uit naar uw reactie.
Met vriendelijke groet,
[Naam bedrijf] - Klantenservice"
$ws.Cells.Item(19, 6).Value = "2025-06-23 18:41:46"
$ws.Cells.Item(19, 7).Value = "Ja"

# --- Extend the conditional-formatting ranges to include row 19 ----------
$ws.Range("D2:D18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D19"))
$ws.Range("G2:G18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G19"))

# --- Dashboard: swap the "Bestelling / Levering" / "Retour / Terugbetaling"
#     summary rows and refresh their counts -------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(3, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(3, 2).Value = 4
$dash.Cells.Item(4, 1).Value = "Bestelling / Levering"
$dash.Cells.Item(4, 2).Value = 3
